$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "Random Forest" (row 7) and "XGBoost" (row 8) rows and insert
# them above row 4, pushing the existing rows 4-8 down to rows 6-10.
$ws.Rows("7:8").Copy()
$ws.Rows("4:5").Insert()

# The inserted rows lose the bordered style on column A; restore it by
# copying the format from a correctly-styled cell in column A.
$ws.Range("A6").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rename the two newly-inserted rows to their "Derived" labels.
$ws.Range("A4").Value = "Derived RF"
$ws.Range("A5").Value = "Derived XGBoost"

# Row 8 ("Multi-Output XGBoost") now carries the values of the XGBoost
# model rather than its own original values; update them accordingly.
$ws.Range("B8").Value = "178.397 ± 6.878"
$ws.Range("C8").Value = "0.928 ± 0.003"
$ws.Range("D8").Value = "984.704 ± 103.751"
$ws.Range("E8").Value = "0.957 ± 0.005"
$ws.Range("F8").Value = "1378.546 ± 96.719"
$ws.Range("G8").Value = "0.958 ± 0.003"
